$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (row 11) down into the
# new rows (12-17) so fonts/number formats match the rest of the table,
# then fill in the actual worklog entries.
$ws.Range("A11:F11").Copy()
$ws.Range("A12:F17").PasteSpecial(-4122)

# Numbers first (log number, date, hours, minutes) - these never touch the
# shared-string table so ordering here is inconsequential.
$ws.Range("A12").Value = 11
$ws.Range("C12").Value = 45562
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 40

$ws.Range("A13").Value = 12
$ws.Range("C13").Value = 45563
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 5

$ws.Range("A14").Value = 13
$ws.Range("C14").Value = 45563
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 0

$ws.Range("A15").Value = 14
$ws.Range("C15").Value = 45563
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 0

$ws.Range("A16").Value = 15
$ws.Range("C16").Value = 45563
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 30

$ws.Range("A17").Value = 16
$ws.Range("C17").Value = 45563
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 0

# Team-member text (reuses names already present in the shared-string
# table from earlier log entries).
$ws.Range("B12").Value = "Andrew"
$ws.Range("B14").Value = "Aiden"
$ws.Range("B15").Value = "Aiden"
$ws.Range("B16").Value = "Aiden"
$ws.Range("B17").Value = "Landon"

# Description text, written in the same order the strings were first
# introduced so the shared-string table layout matches.
$ws.Range("B13").Value = "Landon, Aiden, Andrew"
$ws.Range("F12").Value = "Added all comments to ai shooting behavior code"
$ws.Range("F13").Value = "Discussed ship difficulty logic"
$ws.Range("F15").Value = "Overwrote get_input in AI class to just return coordinate from aiTurn. Added gameloop logic so that you can actually play against AI now. Seems that at least easy mode is working, haven't had a chance to test other modes"
$ws.Range("F14").Value = "Changed AI class to inherit Player so that functions aren't being recreated. Overwrote some functions from Player and GameObject to fit the AI better. Ship placement logic may be done but hasn't been tested and isn't optimized since it's always randomizing without knowledge of previous ship placements or board boundaries and instead error checks each time. Started to modify __take_turn() to take in coordinates from aiTurn but still needs a lot of attention. Waiting for Andrew to push comments on aiTurn function before I start messing with it so that merge conflicts are avoided"
$ws.Range("F16").Value = "Cleaned some more code and implemented medium and hard difficulty for the AI"
$ws.Range("F17").Value = "Made better transitions between Player -> Player and Player -> AI based off of input buffer instead of waiting. Spent a lot of time stress testing with different combinations of AI difficulties, ship amounts, ship directions, and super shot. Updated worklog to be up to date."

# Move the selection (the diff shows the saved view now has B20 selected,
# with no frozen/scrolled topLeftCell override).
$ws.Range("B20").Select()
